$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 13.67700833333333
$ws.Range("N2").Value = 41.031025
$ws.Range("O2").Value = 0.124413831206147
$ws.Range("P2").Value = 0.124413831206147
$ws.Range("Q2").Value = 0.987165430475
$ws.Range("R2").Value = 8.884488874275
$ws.Range("S2").Value = 0.124413831206147
$ws.Range("T2").Value = 0.124413831206147

$ws.Range("M3").Value = 74.64939600000001
$ws.Range("N3").Value = 223.948188
$ws.Range("O3").Value = 0.679053278848249
$ws.Range("P3").Value = 0.6790532788482488
$ws.Range("Q3").Value = 5.387969455092001
$ws.Range("R3").Value = 48.491725095828
$ws.Range("S3").Value = 0.679053278848249
$ws.Range("T3").Value = 0.6790532788482488

$ws.Range("M4").Value = 1.629335666666667
$ws.Range("N4").Value = 4.888007
$ws.Range("O4").Value = 0.01482136207497777
$ws.Range("P4").Value = 0.01482136207497777
$ws.Range("Q4").Value = 0.117600560413
$ws.Range("R4").Value = 1.058405043717
$ws.Range("S4").Value = 0.01482136207497777
$ws.Range("T4").Value = 0.01482136207497777

$ws.Range("M5").Value = 19.17462033333333
$ws.Range("N5").Value = 57.523861
$ws.Range("O5").Value = 0.174423230537864
$ws.Range("P5").Value = 0.174423230537864
$ws.Range("Q5").Value = 1.383966571799
$ws.Range("R5").Value = 12.455699146191
$ws.Range("S5").Value = 0.174423230537864
$ws.Range("T5").Value = 0.174423230537864

$ws.Range("M6").Value = 0.801214
$ws.Range("N6").Value = 2.403642
$ws.Range("O6").Value = 0.007288297332762355
$ws.Range("P6").Value = 0.007288297332762355
$ws.Range("Q6").Value = 0.05782922287800001
$ws.Range("R6").Value = 0.520463005902
$ws.Range("S6").Value = 0.007288297332762355
$ws.Range("T6").Value = 0.007288297332762355
